$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 19.302737645482544
$ws.Range("C2").Value = -4.889486674598956
$ws.Range("D2").Value = -3.6368883229535029
$ws.Range("E2").Value = 0.67431139478685509

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 17.636964166419766
$ws.Range("C3").Value = 2.8609691566184168
$ws.Range("D3").Value = -15.559375505509777
$ws.Range("E3").Value = 5.9981617106704093

# Update the selection to reflect the new active range
$ws.Range("B1:E3").Select()
